$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LP1912")

# Header info: last-updated timestamp + row count
$ws.Range("A2").Value = "Última actualización: 09:23:23"
$ws.Range("A3").Value = "Total filas: 81"

# Data rows 56-86:
#  - rows 56/57 swap places (new arrival data moved them up in the
#    Hora_Llegada-sorted list)
#  - rows 75/76 overwritten with freshly scraped stops
#  - rows 77-86 are newly appended (total rows 72 -> 81)
$rowData = @{}
$rowData[56] = @("08:27:16", "09:17", "27_EL RETIRO", 50, "LP1912")
$rowData[57] = @("07:38:39", "09:17", "14_ABASTO", 99, "LP1912")
$rowData[75] = @("09:23:23", "10:21", "23_HERNANDEZ", 58, "LP1912")
$rowData[76] = @("09:23:23", "10:25", "16_SANTA ANA", 62, "LP1912")
$rowData[77] = @("08:37:25", "10:29", "15_ABASTO", 112, "LP1912")
$rowData[78] = @("09:23:23", "10:29", "14_ABASTO", 66, "LP1912")
$rowData[79] = @("08:45:36", "10:44", "11X44_ETCHEVERRY", 119, "LP1912")
$rowData[80] = @("08:52:50", "10:46", "15_P INDUSTRIAL", 114, "LP1912")
$rowData[81] = @("09:23:23", "10:57", "10_OLMOS", 94, "LP1912")
$rowData[82] = @("09:23:23", "10:59", "27_EL RETIRO", 96, "LP1912")
$rowData[83] = @("09:23:23", "11:01", "81_EL PELIGRO", 98, "LP1912")
$rowData[84] = @("09:23:23", "11:10", "16_P MOR-SANTA ANA", 107, "LP1912")
$rowData[85] = @("09:23:23", "11:14", "14_ABASTO", 111, "LP1912")
$rowData[86] = @("09:23:23", "11:15", "15X38_ABASTO", 112, "LP1912")

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
}

# The other two sheets only get the refreshed "last updated" stamp
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 09:23:23"

$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 09:23:23"

